# "Add Login & Akun Saya"
# Duplicates the existing username/password login row, inserting a new
# account entry above the original one (which shifts two rows down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet starts as:
#   Row1: username | password            (headers)
#   Row2: timothyhudson23@gmail.com (hyperlinked) | IGHOxYmcYdeXBupIKMuP1g==
#
# Target:
#   Row1: username | password
#   Row2: timothyhudson23@gmail.com (hyperlinked - NEW login entry)
#   Row3:           | IGHOxYmcYdeXBupIKMuP1g==
#   Row4: timothyhudson23@gmail.com (hyperlinked - original entry, moved)  | IGHOxYmcYdeXBupIKMuP1g==

# Step 1: push the existing account row down two rows so the original
# record ends up on row 4.
$ws.Rows("2:3").Insert()

# Row insertion doesn't relocate the hyperlink definition that used to sit
# on A2, so drop the stale link before re-creating it on the cell it now
# belongs to.
$ws.Range("A2").Hyperlinks.Delete()

# Step 2: fill in the new login entry (username only on row 2, the
# matching password sits on row 3).
$ws.Range("A2").Value = "timothyhudson23@gmail.com"
$ws.Range("B3").Value = "IGHOxYmcYdeXBupIKMuP1g=="

# Step 3: re-create the hyperlinks - first the original account (now on
# row 4) so it reclaims rId1, then the freshly added one on row 2 as rId2.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:timothyhudson23@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:timothyhudson23@gmail.com")

# Adding a hyperlink re-styles the cell with a fresh xf; put both cells
# back on the shared "Hyperlink" cell style used throughout the sheet.
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A2").Style = "Hyperlink"

# Step 4: match the author's final selection.
$ws.Range("B9").Select()
